$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.619.25'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '1.825.69'
$ws.Range('E3').Value = '  +1.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4674'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3604'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07126'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9018'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07750'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.41'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '1.824.78'
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.267'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.350'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.47'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.009'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008541'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').Value = '26.659.35'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('E22').Value = '  +0.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.904'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.83'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.975'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.862'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08804'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.149'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.829'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.161'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7361'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.440'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01925'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05155'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.876'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5059'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1495'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.041'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('E44').Value = '  +0.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4663'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '97.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.572'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06048'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.79'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.34%  '
